$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ID Kasus"
$ws.Range("B1").Value = "NIK"
$ws.Range("C1").Value = "Nama"
$ws.Range("D1").Value = "Tanggal Lahir"
$ws.Range("E1").Value = "Usia"
$ws.Range("F1").Value = "Jenis Kelamin"
$ws.Range("G1").Value = "Alamat Tempat Tinggal"
$ws.Range("H1").Value = "Kewarganegaraan"
$ws.Range("I1").Value = "No. Telp"
$ws.Range("J1").Value = "Pekerjaan"
$ws.Range("K1").Value = "Status"
$ws.Range("L1").Value = "Tahapan"
$ws.Range("M1").Value = "Hasil"
$ws.Range("N1").Value = "Lokasi saat ini"
$ws.Range("O1").Value = "Tanggal Awal gejala"
$ws.Range("P1").Value = "Gejala"
$ws.Range("Q1").Value = "Riwayat"
$ws.Range("R1").Value = "Tanggal Input"
$ws.Range("S1").Value = "Author"

# Data row
$ws.Range("A2").Value = "covid-1024200001"
$ws.Range("C2").Value = "Liu Xiamei"
$ws.Range("E2").Value = 45
$ws.Range("F2").Value = "L"
$ws.Range("G2").Value = "Mess PT Sansan Melong Asih KOTA CIMAHI Kelurahan Cimahi Selatan Kecamatan Cimahi Selatan"
$ws.Range("H2").Value = "WNA"
$ws.Range("J2").Value = "PT Sansan"
$ws.Range("K2").Value = "ODP"
$ws.Range("L2").Formula = "'1"
$ws.Range("M2").Value = "MENINGGAL"
$ws.Range("N2").Value = "Mess PT Sansan Melong Asih"
$ws.Range("O2").Value = "24/3/2020"
$ws.Range("P2").Formula = "'"
$ws.Range("Q2").Value = "SEMBUH"
$ws.Range("R2").Value = "24/3/2020"
$ws.Range("S2").Value = "Dinkes Kota Cimahi"
